$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
}

# Row 2
$ws.Range("D2").Value = "35.396.34"
$ws.Range("E2").Value = "  +2.80%  "

# Row 3
$ws.Range("D3").Value = "1.847.99"
$ws.Range("E3").Value = "  +2.39%  "

# Row 4
$ws.Range("E4").Value = "  +0.31%  "

# Row 5
Set-TextValue $ws.Range("D5") "229.52"
$ws.Range("E5").Value = "  +1.90%  "

# Row 6
$ws.Range("E6").Value = "  +4.18%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.01"
$ws.Range("E7").Value = "  +0.30%  "

# Row 8
Set-TextValue $ws.Range("D8") "41.97"
$ws.Range("E8").Value = "  +10.34%  "

# Row 9
$ws.Range("E9").Value = "  +6.89%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0691"
$ws.Range("E10").Value = "  +3.41%  "

# Row 11
$ws.Range("E11").Value = "  +3.81%  "

# Row 12
$ws.Range("D12").Value = "2.116.62"
$ws.Range("E12").Value = "  +2.38%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.854.11"
$ws.Range("E13").Value = "  +3.16%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "11.37"
$ws.Range("E14").Value = "  +3.11%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.667"
$ws.Range("E15").Value = "  +6.81%  "

# Row 16
Set-TextValue $ws.Range("D16") "4.66"
$ws.Range("E16").Value = "  +6.65%  "

# Row 17
$ws.Range("D17").Value = "35.410.39"
$ws.Range("E17").Value = "  +2.84%  "

# Row 18
Set-TextValue $ws.Range("D18") "70.22"
$ws.Range("E18").Value = "  +3.70%  "

# Row 19
Set-TextValue $ws.Range("D19") "246.65"
$ws.Range("E19").Value = "  +2.25%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0797"
$ws.Range("E20").Value = "  +4.33%  "

# Row 21
Set-TextValue $ws.Range("D21") "12.04"
$ws.Range("E21").Value = "  +9.03%  "

# Row 22
Set-TextValue $ws.Range("D22") "4.60"
$ws.Range("E22").Value = "  +13.10%  "

# Row 23
Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  +0.25%  "

# Row 24
$ws.Range("E24").Value = "  -0.41%  "

# Row 25
Set-TextValue $ws.Range("D25") "169.06"
$ws.Range("E25").Value = "  -0.59%  "

# Row 26
Set-TextValue $ws.Range("D26") "7.88"
$ws.Range("E26").Value = "  +2.71%  "

# Row 27
Set-TextValue $ws.Range("D27") "17.69"
$ws.Range("E27").Value = "  +1.81%  "

# Row 28
$ws.Range("E28").Value = "  +2.18%  "

# Row 29
$ws.Range("E29").Value = "  +12.97%  "

# Row 30
$ws.Range("E30").Value = "  +0.31%  "

# Row 31
$ws.Range("D31").Value = "3.270.89"
$ws.Range("E31").Value = "  +34.62%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.0542"
$ws.Range("E32").Value = "  +6.12%  "

# Row 33
Set-TextValue $ws.Range("D33") "3.92"
$ws.Range("E33").Value = "  +4.86%  "

# Row 34
Set-TextValue $ws.Range("D34") "4.04"
$ws.Range("E34").Value = "  +6.07%  "

# Row 35
$ws.Range("E35").Value = "  +3.19%  "

# Row 36
Set-TextValue $ws.Range("D36") "97.26"
$ws.Range("E36").Value = "  +19.45%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.681"
$ws.Range("E37").Value = "  +7.27%  "

# Row 38
$ws.Range("D38").Value = "1.351.89"
$ws.Range("E38").Value = "  +1.95%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.08"
$ws.Range("E39").Value = "  +2.83%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.44"
$ws.Range("E40").Value = "  +6.33%  "

# Row 41
$ws.Range("E41").Value = "  +3.13%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.998"
$ws.Range("E42").Value = "  +6.35%  "

# Row 43
$ws.Range("E43").Value = "  +4.01%  "

# Row 44
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D44") "2.48"
$ws.Range("E44").Value = "  +1.16%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D45") "14.62"
$ws.Range("E45").Value = "  +8.20%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.0519"
$ws.Range("E47").Value = "  +1.55%  "

# Row 48
Set-TextValue $ws.Range("D48") "6.17"
$ws.Range("E48").Value = "  +8.14%  "

# Row 49
$ws.Range("D49").Value = "2.015.03"
$ws.Range("E49").Value = "  +2.39%  "

# Row 50
$ws.Range("E50").Value = "  +0.38%  "

# Row 51
Set-TextValue $ws.Range("D51") "103.40"
$ws.Range("E51").Value = "  +1.84%  "
